$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "2016": fill in previously-missing Above Ground (AG) / Below Ground
# (BG) biomass data that had been marked as "missing BG" / "missing AG and
# BG" placeholders.
# ---------------------------------------------------------------------------
$ws16 = $wb.Worksheets.Item("2016")

# Plant-id notes that were not filled in before (column D)
$ws16.Range("D2").Value = "M1E"
$ws16.Range("D3").Value = "M3"
$ws16.Range("D5").Value = "M2"
$ws16.Range("D8").Value = "M1E"
$ws16.Range("D9").Value = "M2"
$ws16.Range("D10").Value = "M2"

# Above-ground (E/F/G) values share one fill-down formula, G2:G7
$ws16.Range("G2:G7").Formula = "=E2-F2"

# Row 7 was missing Below Ground data entirely ("missing BG")
$ws16.Range("H7").Value = 109.1
$ws16.Range("I7").Value = 8.8000000000000007
$ws16.Range("J7").Formula = "=H7-I7"
$ws16.Range("K7").Formula = "=G7/J7"

# Row 8 was missing both AG and BG data ("missing AG and BG")
$ws16.Range("E8").Value = 45.8
$ws16.Range("F8").Value = 14.8
$ws16.Range("G8").Formula = "=E8-F8"
$ws16.Range("H8").Value = 229.9
$ws16.Range("I8").Value = 18.2
$ws16.Range("J8").Formula = "=H8-I8"
$ws16.Range("K8").Formula = "=G8/J8"

# Row 10 was missing Below Ground data ("missing BG")
$ws16.Range("H10").Value = 244.9
$ws16.Range("I10").Value = 15.8
$ws16.Range("J10").Value = 229.1
$ws16.Range("K10").Formula = "=G10/J10"

# Update the Average / Std Dev summary values (K12/K13) to reflect new data
$ws16.Range("K12").Formula = "=AVERAGE(K2:K10)"
$ws16.Range("K13").Formula = "=_xlfn.STDEV.P(K2:K10)"

# The old "missing BG" / "missing AG and BG" notes in column L are no longer
# needed now that the real data has been filled in
[void]$ws16.Range("L7").ClearContents()
[void]$ws16.Range("L8").ClearContents()
[void]$ws16.Range("L10").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "2017": same kind of fix -- fill in previously-missing AG/BG data
# ---------------------------------------------------------------------------
$ws17 = $wb.Worksheets.Item("2017")

# Row 2 was missing AG data entirely ("missing AG")
$ws17.Range("D2").Value = "M2"
$ws17.Range("E2").Value = 86.7
$ws17.Range("F2").Value = 10.1
$ws17.Range("G2").Formula = "=E2-F2"
$ws17.Range("K2").Value = 0.68549000000000004

# Row 5 was missing AG data entirely ("missing AG")
$ws17.Range("E5").Value = 39.6
$ws17.Range("F5").Value = 22.3
$ws17.Range("G5").Formula = "=E5-F5"
$ws17.Range("K5").Value = 0.55897459999999999

# Row 7 was missing AG data entirely ("missing AG")
$ws17.Range("E7").Value = 34.1
$ws17.Range("F7").Value = 19.2
$ws17.Range("G7").Formula = "=E7-F7"
$ws17.Range("K7").Formula = "=G7/J7"

# Row 8 was missing AG data entirely ("missing AG")
$ws17.Range("E8").Value = 59.3
$ws17.Range("F8").Value = 21.7
$ws17.Range("G8").Formula = "=E8-F8"
$ws17.Range("K8").Formula = "=G8/J8"

# K9 recomputed explicitly (was part of the shared K3:K10 formula group)
$ws17.Range("K9").Formula = "=G9/J9"

# Update the Average / Std Dev summary values (K13/K14) to reflect new data
$ws17.Range("K13").Formula = "=AVERAGE(K2:K10)"
$ws17.Range("K14").Formula = "=_xlfn.STDEV.P(K2:K10)"

# The old "missing AG" notes in column L are no longer needed now that the
# real data has been filled in
[void]$ws17.Range("L2").ClearContents()
[void]$ws17.Range("L5").ClearContents()
[void]$ws17.Range("L7").ClearContents()
[void]$ws17.Range("L8").ClearContents()
